# Auto-applied edit: update "Error" sheet per commit diff (rows 2-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Set cell values ---
$ws.Range("B2").Value = 45972
$ws.Range("A3").Value = 251310
$ws.Range("B3").Value = 45770
$ws.Range("D3").Value = 45769
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 8611
$ws.Range("G3").Value = 336
$ws.Range("H3").Value = "6"
$ws.Range("J3").Value = 410
$ws.Range("K3").Value = 820
$ws.Range("A4").Value = 252459
$ws.Range("B4").Value = 45855
$ws.Range("D4").Value = "CAMPO VUOTO"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 4000
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = "CAMPO VUOTO"
$ws.Range("I4").Value = "foglio"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 480
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = 3
$ws.Range("A5").Value = "CAMPO VUOTO"
$ws.Range("B5").Value = "CAMPO VUOTO"
$ws.Range("C5").Value = "DA STAMPARE"
$ws.Range("D5").Value = "CAMPO VUOTO"
$ws.Range("E5").Value = "CAMPO VUOTO"
$ws.Range("F5").Value = "CAMPO VUOTO"
$ws.Range("G5").Value = "CAMPO VUOTO"
$ws.Range("H5").Value = "CAMPO VUOTO"
$ws.Range("I5").Value = "foglio"
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "CAMPO VUOTO"
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = "CAMPO VUOTO"
$ws.Range("N5").Value = "CAMPO VUOTO"
$ws.Range("O5").Value = "CAMPO VUOTO"
$ws.Range("P5").Value = "CAMPO VUOTO"
$ws.Range("A6").Value = "CAMPO VUOTO"
$ws.Range("B6").Value = "CAMPO VUOTO"
$ws.Range("C6").Value = "DA STAMPARE"
$ws.Range("E6").Value = "CAMPO VUOTO"
$ws.Range("F6").Value = "CAMPO VUOTO"
$ws.Range("G6").Value = "CAMPO VUOTO"
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = "CAMPO VUOTO"
$ws.Range("M6").Value = "CAMPO VUOTO"
$ws.Range("N6").Value = "CAMPO VUOTO"
$ws.Range("A7").Value = "CAMPO VUOTO"
$ws.Range("B7").Value = "CAMPO VUOTO"
$ws.Range("C7").Value = "DA STAMPARE"
$ws.Range("D7").Value = "CAMPO VUOTO"
$ws.Range("E7").Value = "CAMPO VUOTO"
$ws.Range("F7").Value = "CAMPO VUOTO"
$ws.Range("G7").Value = "CAMPO VUOTO"
$ws.Range("H7").Value = "CAMPO VUOTO"
$ws.Range("I7").Value = "foglio"
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = "CAMPO VUOTO"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = "CAMPO VUOTO"
$ws.Range("N7").Value = "CAMPO VUOTO"
$ws.Range("O7").Value = "CAMPO VUOTO"
$ws.Range("P7").Value = "CAMPO VUOTO"
$ws.Range("A8").Value = 253258
$ws.Range("B8").Value = 45922
$ws.Range("C8").Value = "STAMPATO"
$ws.Range("D8").Value = "CAMPO VUOTO"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 4805
$ws.Range("G8").Value = 223
$ws.Range("H8").Value = "2"
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 76
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = "CAMPO VUOTO"
$ws.Range("P8").Value = "CAMPO VUOTO"
$ws.Range("A9").Value = "CAMPO VUOTO"
$ws.Range("B9").Value = "CAMPO VUOTO"
$ws.Range("C9").Value = "DA STAMPARE"
$ws.Range("D9").Value = "CAMPO VUOTO"
$ws.Range("E9").Value = "CAMPO VUOTO"
$ws.Range("F9").Value = "CAMPO VUOTO"
$ws.Range("G9").Value = "CAMPO VUOTO"
$ws.Range("H9").Value = "CAMPO VUOTO"
$ws.Range("I9").Value = "foglio"
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = "CAMPO VUOTO"
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = "CAMPO VUOTO"
$ws.Range("N9").Value = "CAMPO VUOTO"
$ws.Range("O9").Value = "CAMPO VUOTO"
$ws.Range("P9").Value = "CAMPO VUOTO"
$ws.Range("A10").Value = 254225
$ws.Range("B10").Value = 45975
$ws.Range("C10").Value = "IN STAMPA"
$ws.Range("D10").Value = 45960
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 3943
$ws.Range("G10").Value = 200
$ws.Range("H10").Value = "CAMPO VUOTO"
$ws.Range("I10").Value = "bobina"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 780
$ws.Range("L10").Value = "CAMPO VUOTO"
$ws.Range("M10").Value = "Dati OK"
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = "X"
$ws.Range("P10").Value = 40635

# --- 2) Fix up cell styles that changed (reuse existing style slots) ---
# Style reference cells that keep their style identity throughout this edit:
#   A2 -> default/general style (style index 0)
#   B2 -> date style "yyyy-mm-dd h:mm:ss" (style index 1)
#   D2 -> "CAMPO VUOTO" yellow-fill style (style index 2)

$ws.Range("A2").Copy()
foreach ($addr in @("L4","L5","L7","L8")) { $ws.Range($addr).PasteSpecial(-4122) }

$ws.Range("B2").Copy()
foreach ($addr in @("B10","D10")) { $ws.Range($addr).PasteSpecial(-4122) }

$ws.Range("D2").Copy()
foreach ($addr in @("D4","A5","B5","D5","E5","F5","G5","K5","M5","N5","A6","B6","E6","F6","G6","K6","M6","N6","A7","B7","D7","E7","F7","G7","K7","M7","N7","D8","A9","B9","D9","E9","F9","G9","K9","M9","N9","L10")) { $ws.Range($addr).PasteSpecial(-4122) }

$excel.CutCopyMode = 0
